# "fall 22 week 13 day-after inputs"
#
# The "Games Remaining Needed" column (T) could go negative once a team had
# already clinched / fallen out of reach, which doesn't make sense for a
# countdown of games still needed. Wrap the existing formula in MAX(...,0)
# for every team row in both week blocks (rows 3-10 and 15-22) so it floors
# at zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wookie Mistakes")

$rows = @(3,4,5,6,7,8,9,10,15,16,17,18,19,20,21,22)
foreach ($r in $rows) {
    $ws.Range("T$r").Formula = '=MAX(S' + $r + '-COUNTIF(B' + $r + ':Q' + $r + ', "W")-COUNTIF(B' + $r + ':Q' + $r + ', "L"), 0)'
}

# Move the selection to where the author's cursor ended up after making the
# edits (also nudges the scrolled-into-view region toward column O).
$ws.Range("O1").Select()
$ws.Range("U3").Select()
